# Update SRA_data sheet:
#  - Column J (design_description): "Sequencing performed at None"
#    -> "Sequencing performed at Oregon State University Center for Quantitative Life Sciences Genomics Core"
#  - Column L (filename): derived from column B (library_ID) + "_R1.fastq.gz"
#  - Column M (filename2): derived from column B (library_ID) + "_R2.fastq.gz"
# for all data rows (2 through 107) on the "SRA_data" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SRA_data")

$newSeqText = "Sequencing performed at Oregon State University Center for Quantitative Life Sciences Genomics Core"

$firstRow = 2
$lastRow = 107

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $libraryId = $ws.Cells.Item($row, 2).Value2

    if ($libraryId -ne $null -and $libraryId -ne "") {
        $ws.Cells.Item($row, 10).Value = $newSeqText
        $ws.Cells.Item($row, 12).Value = $libraryId + "_R1.fastq.gz"
        $ws.Cells.Item($row, 13).Value = $libraryId + "_R2.fastq.gz"
    }
}
